$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 2 de Agosto de 2020 a las 06:56"

# Refresh case counts for Peru (row 10)
$ws.Cells.Item(10, 2).Value = 422183
$ws.Cells.Item(10, 4).Value = 290835
$ws.Cells.Item(10, 5).Value = 111940
$ws.Cells.Item(10, 8).Value = 19408

# Refresh case counts for Pakistan (row 16)
$ws.Cells.Item(16, 2).Value = 279146
$ws.Cells.Item(16, 4).Value = 248027
$ws.Cells.Item(16, 5).Value = 25149
$ws.Cells.Item(16, 8).Value = 5970

# Refresh case counts for Australia (row 72)
$ws.Cells.Item(72, 2).Value = 17895
$ws.Cells.Item(72, 3).Value = 613
$ws.Cells.Item(72, 4).Value = 10204
$ws.Cells.Item(72, 5).Value = 7483
$ws.Cells.Item(72, 7).Value = 7
$ws.Cells.Item(72, 8).Value = 208

# Haiti's updated case counts push it above Tayikistan and Finlandia in the
# ranking (sorted descending by total cases), so rows 91-93 are re-sorted:
# row 91 becomes Haiti (with its fresh numbers), while Tayikistan and
# Finlandia shift down one row each, keeping their own existing numbers.
$ws.Cells.Item(91, 1).Value = "Haiti"
$ws.Cells.Item(91, 2).Value = 7468
$ws.Cells.Item(91, 3).Value = 44
$ws.Cells.Item(91, 4).Value = 4606
$ws.Cells.Item(91, 5).Value = 2697
$ws.Cells.Item(91, 6).Value = 0
$ws.Cells.Item(91, 7).Value = 4
$ws.Cells.Item(91, 8).Value = 165

$ws.Cells.Item(92, 1).Value = "Tayikistan"
$ws.Cells.Item(92, 2).Value = 7451
$ws.Cells.Item(92, 3).Value = 0
$ws.Cells.Item(92, 4).Value = 6233
$ws.Cells.Item(92, 5).Value = 1158
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 60

$ws.Cells.Item(93, 1).Value = "Finlandia"
$ws.Cells.Item(93, 2).Value = 7443
$ws.Cells.Item(93, 3).Value = 0
$ws.Cells.Item(93, 4).Value = 6950
$ws.Cells.Item(93, 5).Value = 164
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 0
$ws.Cells.Item(93, 8).Value = 329

# Refresh case counts for Butan (row 189)
$ws.Cells.Item(189, 2).Value = 102
$ws.Cells.Item(189, 3).Value = 1
$ws.Cells.Item(189, 5).Value = 13
